$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A10").NumberFormat = "0.00E+00"
Write-Host "Done"
